# Applies the cryptos.xlsx price/volume(1h) refresh described in the commit
# "Updated symbol list on Thu Feb 16 23:42:38 UTC 2023 with GitHub Actions".
# Columns D (Price) and E (Volume(1h)) are stored as literal text, so each
# write uses a leading apostrophe to force text-entry (matching the original
# inlineStr cells) and then resets Style to "Normal" so no stray number format
# (e.g. quote-prefix / @ text format) gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $CellRef, $Text) {
    $range = $Sheet.Range($CellRef)
    $range.Value = "'" + $Text
    $range.Style = "Normal"
}

Set-TextValue $ws "D2" "308.85"
Set-TextValue $ws "E2" "-2.28%"
Set-TextValue $ws "D3" "48.27"
Set-TextValue $ws "E3" "6.61%"
Set-TextValue $ws "D4" "5.199"
Set-TextValue $ws "E4" "0.43%"
Set-TextValue $ws "D5" "0.07726"
Set-TextValue $ws "E5" "-4.41%"
Set-TextValue $ws "D6" "4.507"
Set-TextValue $ws "E6" "-0.62%"
Set-TextValue $ws "D7" "1.296"
Set-TextValue $ws "E7" "18.58%"
Set-TextValue $ws "E8" "-6.82%"
Set-TextValue $ws "D9" "0.1230"
Set-TextValue $ws "E9" "-5.82%"
Set-TextValue $ws "D10" "0.1921"
Set-TextValue $ws "E10" "-0.76%"
Set-TextValue $ws "D11" "0.09182"
Set-TextValue $ws "E11" "-3.38%"
Set-TextValue $ws "D12" "0.04551"
Set-TextValue $ws "E12" "7.48%"
Set-TextValue $ws "D13" "0.1049"
Set-TextValue $ws "E13" "0.47%"
Set-TextValue $ws "D14" "0.001291"
Set-TextValue $ws "E14" "-1.81%"
Set-TextValue $ws "D15" "0.04205"
Set-TextValue $ws "E15" "-1.81%"
Set-TextValue $ws "D16" "0.005870"
Set-TextValue $ws "E16" "-1.14%"
Set-TextValue $ws "D17" "3.343"
Set-TextValue $ws "E17" "-1.69%"
Set-TextValue $ws "D18" "2.402"
Set-TextValue $ws "E18" "-0.34%"
Set-TextValue $ws "E19" "2.07%"
Set-TextValue $ws "D20" "8.106"
Set-TextValue $ws "E20" "-1.24%"
Set-TextValue $ws "D21" "0.1374"
Set-TextValue $ws "E21" "-0.75%"
Set-TextValue $ws "D22" "0.3033"
Set-TextValue $ws "E22" "-3.56%"
Set-TextValue $ws "E23" "1.45%"
Set-TextValue $ws "D24" "0.004092"
Set-TextValue $ws "E24" "-3.16%"
Set-TextValue $ws "E25" "1.14%"
Set-TextValue $ws "D26" "0.0003566"
Set-TextValue $ws "E26" "-95.19%"
Set-TextValue $ws "D38" "0.02562"
Set-TextValue $ws "E38" "-5.49%"
Set-TextValue $ws "D39" "0.05720"
Set-TextValue $ws "E39" "4.66%"
Set-TextValue $ws "D40" "0.01085"
Set-TextValue $ws "E40" "84.93%"
Set-TextValue $ws "D41" "0.007965"
Set-TextValue $ws "E41" "2.40%"
Set-TextValue $ws "E42" "-0.39%"
Set-TextValue $ws "D43" "0.008392"
Set-TextValue $ws "D44" "0.007780"
Set-TextValue $ws "E44" "-9.51%"
Set-TextValue $ws "E45" "7.66%"
Set-TextValue $ws "D46" "0.00006842"
Set-TextValue $ws "E46" "0.58%"
Set-TextValue $ws "E47" "1.04%"
Set-TextValue $ws "D48" "0.05552"
Set-TextValue $ws "E48" "-10.84%"
Set-TextValue $ws "D49" "0.004030"
Set-TextValue $ws "E49" "1.13%"
Set-TextValue $ws "D50" "0.00002116"
Set-TextValue $ws "E50" "1.04%"
Set-TextValue $ws "D51" "0.0002015"
Set-TextValue $ws "E51" "1.04%"
